# Refresh generated schema files: add DiseaseListEntry and MatrixDiseaseList sheets
$wb = $excel.ActiveWorkbook

# --- Add "DiseaseListEntry" sheet (header row with many schema columns) ---
$wsEntry = $wb.Worksheets.Add()
$wsEntry.Name = "DiseaseListEntry"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsEntry.Move($null, $lastSheet)
# Re-fetch the worksheet reference after the move so subsequent writes
# land on the correct (moved) sheet.
$wsEntry = $wb.Worksheets.Item("DiseaseListEntry")

$entryHeaders = @(
    "category_class",
    "label",
    "definition",
    "synonyms",
    "subsets",
    "crossreferences",
    "is_matrix_manually_excluded",
    "is_matrix_manually_included",
    "is_clingen",
    "is_grouping_subset",
    "is_grouping_subset_ancestor",
    "is_orphanet_subtype",
    "is_orphanet_subtype_descendant",
    "is_omimps",
    "is_omimps_descendant",
    "is_leaf",
    "is_leaf_direct_parent",
    "is_orphanet_disorder",
    "is_omim",
    "is_icd_category",
    "is_icd_chapter_code",
    "is_icd_chapter_header",
    "is_icd_billable",
    "is_mondo_subtype",
    "is_pathway_defect",
    "is_susceptibility",
    "is_paraphilic",
    "is_acquired",
    "is_andor",
    "is_withorwithout",
    "is_obsoletion_candidate",
    "is_unclassified_hereditary",
    "official_matrix_filter",
    "harrisons_view",
    "mondo_txgnn",
    "mondo_top_grouping",
    "medical_specialization",
    "txgnn",
    "anatomical",
    "is_pathogen_caused",
    "is_cancer",
    "is_glucose_dysfunction",
    "tag_existing_treatment",
    "tag_qaly_lost",
    "subset_group_id",
    "subset_group_label",
    "other_subsets_count"
)

for ($i = 0; $i -lt $entryHeaders.Length; $i++) {
    $wsEntry.Cells.Item(1, $i + 1).Value = $entryHeaders[$i]
}

# --- Add "MatrixDiseaseList" sheet (single "entries" header) ---
$wsList = $wb.Worksheets.Add()
$wsList.Name = "MatrixDiseaseList"
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsList.Move($null, $lastSheet2)
$wsList = $wb.Worksheets.Item("MatrixDiseaseList")

$wsList.Cells.Item(1, 1).Value = "entries"
